# "add functions to settings and name manager after inserting demo sheet"
#
# The demo sheet was inserted while the workbook-settings custom functions
# (ADD_TWO_DAYS, CALCULATE_AREA, IN_RANGE, JOIN_STRINGS, TO_POWER) had not
# yet been (re)registered in the Name Manager, so the demo formula that
# calls one of them resolves to #NAME?. Remove the defined names that
# back those functions.

$wb = $excel.ActiveWorkbook

# Collect the defined-name collection into a plain array first -- deleting
# while iterating the live COM collection would skip entries.
$namesToRemove = @()
foreach ($n in $wb.Names) {
    $namesToRemove += $n.Name
}
foreach ($name in $namesToRemove) {
    $wb.Names($name).Delete()
}

# Force the demo sheet's formula to re-resolve against the now-missing
# names so the cached result reflects the #NAME? error instead of the
# stale value that was cached while the names still existed.
$ws = $wb.Worksheets.Item("DemoFunctions")
$cell = $ws.Range("B4")
$cell.Formula = $cell.Formula
